# Fruta / hortaliza, semanal
# Inserts two new weekly price records (rows 502 and 503) into the "Uva"
# (grape) sheet of the Macroferia Regional de Talca data set. Inserting the
# rows shifts every subsequent record down by two rows (old row 502 becomes
# row 504, ..., old row 539 becomes row 541), which also grows the table by
# two rows overall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 502, pushing all the
# following rows (and the sheet dimension) down by two.
$ws.Rows("502:503").Insert()

# ---- New row 502: Flame Seedless, Provincia de Limarí -----------------
$ws.Range("A502").Value2 = 5
$ws.Range("B502").Value2 = "Macroferia Regional de Talca"
$ws.Range("C502").Value2 = "Maule"
$ws.Range("D502").Value2 = 44931
$ws.Range("E502").Value2 = 7
$ws.Range("F502").Value2 = "Fruta"
$ws.Range("G502").Value2 = 100109
$ws.Range("H502").Value2 = "Uva"
$ws.Range("I502").Value2 = 100109001
$ws.Range("J502").Value2 = "Uva"
$ws.Range("K502").Value2 = "Flame Seedless"
$ws.Range("L502").Value2 = "Primera"
$ws.Range("M502").Value2 = 150
$ws.Range("N502").Value2 = 12000
$ws.Range("O502").Value2 = 12000
$ws.Range("P502").Value2 = 12000
$ws.Range("Q502").Value2 = "$/bandeja 12 kilos"
$ws.Range("R502").Value2 = "Provincia de Limarí"
$ws.Range("S502").Value2 = 1000
$ws.Range("T502").Value2 = 12

# ---- New row 503: Superior Seedless, Provincia de Limarí --------------
$ws.Range("A503").Value2 = 5
$ws.Range("B503").Value2 = "Macroferia Regional de Talca"
$ws.Range("C503").Value2 = "Maule"
$ws.Range("D503").Value2 = 44931
$ws.Range("E503").Value2 = 7
$ws.Range("F503").Value2 = "Fruta"
$ws.Range("G503").Value2 = 100109
$ws.Range("H503").Value2 = "Uva"
$ws.Range("I503").Value2 = 100109001
$ws.Range("J503").Value2 = "Uva"
$ws.Range("K503").Value2 = "Superior Seedless"
$ws.Range("L503").Value2 = "Primera"
$ws.Range("M503").Value2 = 80
$ws.Range("N503").Value2 = 14000
$ws.Range("O503").Value2 = 14000
$ws.Range("P503").Value2 = 14000
$ws.Range("Q503").Value2 = "$/bandeja 10 kilos"
$ws.Range("R503").Value2 = "Provincia de Limarí"
$ws.Range("S503").Value2 = 1400
$ws.Range("T503").Value2 = 10
